$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Cxcl12"
$ws.Cells.Item(2,3).Value = "Itga4"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = [double]"3"
$ws.Cells.Item(2,6).Value = [double]"1"
$ws.Cells.Item(2,7).Value = [double]"81.05837566666666"
$ws.Cells.Item(2,8).Value = [double]"243.175127"
$ws.Cells.Item(2,9).Value = [double]"0.3545816884225585"
$ws.Cells.Item(2,10).Value = [double]"0.3545816884225585"
$ws.Cells.Item(2,11).Value = [double]"2"
$ws.Cells.Item(2,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(2,13).Value = [double]"23.65990166666667"
$ws.Cells.Item(2,14).Value = [double]"70.979705"
$ws.Cells.Item(2,15).Value = [double]"0.2997993941754699"
$ws.Cells.Item(2,16).Value = [double]"0.29979939417547"
$ws.Cells.Item(2,17).Value = [double]"1917.833197533059"
$ws.Cells.Item(2,18).Value = [double]"17260.49877779753"
$ws.Cells.Item(2,19).Value = [double]"0.1063033753747983"
$ws.Cells.Item(2,20).Value = [double]"0.1063033753747983"

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Cxcl12"
$ws.Cells.Item(3,3).Value = "Itga4"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = [double]"3"
$ws.Cells.Item(3,6).Value = [double]"1"
$ws.Cells.Item(3,7).Value = [double]"81.05837566666666"
$ws.Cells.Item(3,8).Value = [double]"243.175127"
$ws.Cells.Item(3,9).Value = [double]"0.3545816884225585"
$ws.Cells.Item(3,10).Value = [double]"0.3545816884225585"
$ws.Cells.Item(3,11).Value = [double]"1"
$ws.Cells.Item(3,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(3,13).Value = [double]"0.07690566666666666"
$ws.Cells.Item(3,14).Value = [double]"0.230717"
$ws.Cells.Item(3,15).Value = [double]"0.0009744872400636476"
$ws.Cells.Item(3,16).Value = [double]"0.0009744872400636479"
$ws.Cells.Item(3,17).Value = [double]"6.23384841956211"
$ws.Cells.Item(3,18).Value = [double]"56.104635776059"
$ws.Cells.Item(3,19).Value = [double]"0.0003455353309280072"
$ws.Cells.Item(3,20).Value = [double]"0.0003455353309280074"

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Cxcl12"
$ws.Cells.Item(4,3).Value = "Itga4"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = [double]"3"
$ws.Cells.Item(4,6).Value = [double]"1"
$ws.Cells.Item(4,7).Value = [double]"81.05837566666666"
$ws.Cells.Item(4,8).Value = [double]"243.175127"
$ws.Cells.Item(4,9).Value = [double]"0.3545816884225585"
$ws.Cells.Item(4,10).Value = [double]"0.3545816884225585"
$ws.Cells.Item(4,11).Value = [double]"3"
$ws.Cells.Item(4,12).Value = [double]"1"
$ws.Cells.Item(4,13).Value = [double]"53.21452433333334"
$ws.Cells.Item(4,14).Value = [double]"159.643573"
$ws.Cells.Item(4,15).Value = [double]"0.6742919890890982"
$ws.Cells.Item(4,16).Value = [double]"0.6742919890890983"
$ws.Cells.Item(4,17).Value = [double]"4313.482904334308"
$ws.Cells.Item(4,18).Value = [double]"38821.34613900877"
$ws.Cells.Item(4,19).Value = [double]"0.2390915919810178"
$ws.Cells.Item(4,20).Value = [double]"0.2390915919810178"

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Cxcl12"
$ws.Cells.Item(5,3).Value = "Itga4"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = [double]"3"
$ws.Cells.Item(5,6).Value = [double]"1"
$ws.Cells.Item(5,7).Value = [double]"81.05837566666666"
$ws.Cells.Item(5,8).Value = [double]"243.175127"
$ws.Cells.Item(5,9).Value = [double]"0.3545816884225585"
$ws.Cells.Item(5,10).Value = [double]"0.3545816884225585"
$ws.Cells.Item(5,11).Value = [double]"3"
$ws.Cells.Item(5,12).Value = [double]"1"
$ws.Cells.Item(5,13).Value = [double]"1.967779333333333"
$ws.Cells.Item(5,14).Value = [double]"5.903338"
$ws.Cells.Item(5,15).Value = [double]"0.02493412949536815"
$ws.Cells.Item(5,16).Value = [double]"0.02493412949536816"
$ws.Cells.Item(5,17).Value = [double]"159.5049964304362"
$ws.Cells.Item(5,18).Value = [double]"1435.544967873926"
$ws.Cells.Item(5,19).Value = [double]"0.008841185735814355"
$ws.Cells.Item(5,20).Value = [double]"0.008841185735814357"

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Cxcl12"
$ws.Cells.Item(6,3).Value = "Itga4"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = [double]"3"
$ws.Cells.Item(6,6).Value = [double]"1"
$ws.Cells.Item(6,7).Value = [double]"123.018252"
$ws.Cells.Item(6,8).Value = [double]"369.054756"
$ws.Cells.Item(6,9).Value = [double]"0.5381309351710768"
$ws.Cells.Item(6,10).Value = [double]"0.5381309351710768"
$ws.Cells.Item(6,11).Value = [double]"2"
$ws.Cells.Item(6,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(6,13).Value = [double]"23.65990166666667"
$ws.Cells.Item(6,14).Value = [double]"70.979705"
$ws.Cells.Item(6,15).Value = [double]"0.2997993941754699"
$ws.Cells.Item(6,16).Value = [double]"0.29979939417547"
$ws.Cells.Item(6,17).Value = [double]"2910.59974552522"
$ws.Cells.Item(6,18).Value = [double]"26195.39770972698"
$ws.Cells.Item(6,19).Value = [double]"0.1613313283513679"
$ws.Cells.Item(6,20).Value = [double]"0.1613313283513679"

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Cxcl12"
$ws.Cells.Item(7,3).Value = "Itga4"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = [double]"3"
$ws.Cells.Item(7,6).Value = [double]"1"
$ws.Cells.Item(7,7).Value = [double]"123.018252"
$ws.Cells.Item(7,8).Value = [double]"369.054756"
$ws.Cells.Item(7,9).Value = [double]"0.5381309351710768"
$ws.Cells.Item(7,10).Value = [double]"0.5381309351710768"
$ws.Cells.Item(7,11).Value = [double]"1"
$ws.Cells.Item(7,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(7,13).Value = [double]"0.07690566666666666"
$ws.Cells.Item(7,14).Value = [double]"0.230717"
$ws.Cells.Item(7,15).Value = [double]"0.0009744872400636476"
$ws.Cells.Item(7,16).Value = [double]"0.0009744872400636479"
$ws.Cells.Item(7,17).Value = [double]"9.460800682227999"
$ws.Cells.Item(7,18).Value = [double]"85.14720614005201"
$ws.Cells.Item(7,19).Value = [double]"0.0005244017298077323"
$ws.Cells.Item(7,20).Value = [double]"0.0005244017298077324"

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Cxcl12"
$ws.Cells.Item(8,3).Value = "Itga4"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = [double]"3"
$ws.Cells.Item(8,6).Value = [double]"1"
$ws.Cells.Item(8,7).Value = [double]"123.018252"
$ws.Cells.Item(8,8).Value = [double]"369.054756"
$ws.Cells.Item(8,9).Value = [double]"0.5381309351710768"
$ws.Cells.Item(8,10).Value = [double]"0.5381309351710768"
$ws.Cells.Item(8,11).Value = [double]"3"
$ws.Cells.Item(8,12).Value = [double]"1"
$ws.Cells.Item(8,13).Value = [double]"53.21452433333334"
$ws.Cells.Item(8,14).Value = [double]"159.643573"
$ws.Cells.Item(8,15).Value = [double]"0.6742919890890982"
$ws.Cells.Item(8,16).Value = [double]"0.6742919890890983"
$ws.Cells.Item(8,17).Value = [double]"6546.357764498132"
$ws.Cells.Item(8,18).Value = [double]"58917.21988048319"
$ws.Cells.Item(8,19).Value = [double]"0.3628573786668819"
$ws.Cells.Item(8,20).Value = [double]"0.3628573786668819"

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Cxcl12"
$ws.Cells.Item(9,3).Value = "Itga4"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = [double]"3"
$ws.Cells.Item(9,6).Value = [double]"1"
$ws.Cells.Item(9,7).Value = [double]"123.018252"
$ws.Cells.Item(9,8).Value = [double]"369.054756"
$ws.Cells.Item(9,9).Value = [double]"0.5381309351710768"
$ws.Cells.Item(9,10).Value = [double]"0.5381309351710768"
$ws.Cells.Item(9,11).Value = [double]"3"
$ws.Cells.Item(9,12).Value = [double]"1"
$ws.Cells.Item(9,13).Value = [double]"1.967779333333333"
$ws.Cells.Item(9,14).Value = [double]"5.903338"
$ws.Cells.Item(9,15).Value = [double]"0.02493412949536815"
$ws.Cells.Item(9,16).Value = [double]"0.02493412949536816"
$ws.Cells.Item(9,17).Value = [double]"242.072773908392"
$ws.Cells.Item(9,18).Value = [double]"2178.654965175528"
$ws.Cells.Item(9,19).Value = [double]"0.01341782642301919"
$ws.Cells.Item(9,20).Value = [double]"0.0134178264230192"

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Cxcl12"
$ws.Cells.Item(10,3).Value = "Itga4"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = [double]"3"
$ws.Cells.Item(10,6).Value = [double]"1"
$ws.Cells.Item(10,7).Value = [double]"0.3624666666666667"
$ws.Cells.Item(10,8).Value = [double]"1.0874"
$ws.Cells.Item(10,9).Value = [double]"0.001585573873230423"
$ws.Cells.Item(10,10).Value = [double]"0.001585573873230423"
$ws.Cells.Item(10,11).Value = [double]"2"
$ws.Cells.Item(10,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(10,13).Value = [double]"23.65990166666667"
$ws.Cells.Item(10,14).Value = [double]"70.979705"
$ws.Cells.Item(10,15).Value = [double]"0.2997993941754699"
$ws.Cells.Item(10,16).Value = [double]"0.29979939417547"
$ws.Cells.Item(10,17).Value = [double]"8.575925690777778"
$ws.Cells.Item(10,18).Value = [double]"77.18333121699999"
$ws.Cells.Item(10,19).Value = [double]"0.0004753540866149342"
$ws.Cells.Item(10,20).Value = [double]"0.0004753540866149343"

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Cxcl12"
$ws.Cells.Item(11,3).Value = "Itga4"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = [double]"3"
$ws.Cells.Item(11,6).Value = [double]"1"
$ws.Cells.Item(11,7).Value = [double]"0.3624666666666667"
$ws.Cells.Item(11,8).Value = [double]"1.0874"
$ws.Cells.Item(11,9).Value = [double]"0.001585573873230423"
$ws.Cells.Item(11,10).Value = [double]"0.001585573873230423"
$ws.Cells.Item(11,11).Value = [double]"1"
$ws.Cells.Item(11,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(11,13).Value = [double]"0.07690566666666666"
$ws.Cells.Item(11,14).Value = [double]"0.230717"
$ws.Cells.Item(11,15).Value = [double]"0.0009744872400636476"
$ws.Cells.Item(11,16).Value = [double]"0.0009744872400636479"
$ws.Cells.Item(11,17).Value = [double]"0.02787574064444444"
$ws.Cells.Item(11,18).Value = [double]"0.2508816658"
$ws.Cells.Item(11,19).Value = [double]"1.545121507641343E-06"
$ws.Cells.Item(11,20).Value = [double]"1.545121507641344E-06"

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Cxcl12"
$ws.Cells.Item(12,3).Value = "Itga4"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = [double]"3"
$ws.Cells.Item(12,6).Value = [double]"1"
$ws.Cells.Item(12,7).Value = [double]"0.3624666666666667"
$ws.Cells.Item(12,8).Value = [double]"1.0874"
$ws.Cells.Item(12,9).Value = [double]"0.001585573873230423"
$ws.Cells.Item(12,10).Value = [double]"0.001585573873230423"
$ws.Cells.Item(12,11).Value = [double]"3"
$ws.Cells.Item(12,12).Value = [double]"1"
$ws.Cells.Item(12,13).Value = [double]"53.21452433333334"
$ws.Cells.Item(12,14).Value = [double]"159.643573"
$ws.Cells.Item(12,15).Value = [double]"0.6742919890890982"
$ws.Cells.Item(12,16).Value = [double]"0.6742919890890983"
$ws.Cells.Item(12,17).Value = [double]"19.28849125335556"
$ws.Cells.Item(12,18).Value = [double]"173.5964212802"
$ws.Cells.Item(12,19).Value = [double]"0.001069139760828248"
$ws.Cells.Item(12,20).Value = [double]"0.001069139760828248"

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Cxcl12"
$ws.Cells.Item(13,3).Value = "Itga4"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = [double]"3"
$ws.Cells.Item(13,6).Value = [double]"1"
$ws.Cells.Item(13,7).Value = [double]"0.3624666666666667"
$ws.Cells.Item(13,8).Value = [double]"1.0874"
$ws.Cells.Item(13,9).Value = [double]"0.001585573873230423"
$ws.Cells.Item(13,10).Value = [double]"0.001585573873230423"
$ws.Cells.Item(13,11).Value = [double]"3"
$ws.Cells.Item(13,12).Value = [double]"1"
$ws.Cells.Item(13,13).Value = [double]"1.967779333333333"
$ws.Cells.Item(13,14).Value = [double]"5.903338"
$ws.Cells.Item(13,15).Value = [double]"0.02493412949536815"
$ws.Cells.Item(13,16).Value = [double]"0.02493412949536816"
$ws.Cells.Item(13,17).Value = [double]"0.7132544156888889"
$ws.Cells.Item(13,18).Value = [double]"6.419289741199999"
$ws.Cells.Item(13,19).Value = [double]"3.953490427959981E-05"
$ws.Cells.Item(13,20).Value = [double]"3.953490427959983E-05"

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Cxcl12"
$ws.Cells.Item(14,3).Value = "Itga4"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = [double]"3"
$ws.Cells.Item(14,6).Value = [double]"1"
$ws.Cells.Item(14,7).Value = [double]"24.16373066666667"
$ws.Cells.Item(14,8).Value = [double]"72.491192"
$ws.Cells.Item(14,9).Value = [double]"0.1057018025331343"
$ws.Cells.Item(14,10).Value = [double]"0.1057018025331344"
$ws.Cells.Item(14,11).Value = [double]"2"
$ws.Cells.Item(14,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(14,13).Value = [double]"23.65990166666667"
$ws.Cells.Item(14,14).Value = [double]"70.979705"
$ws.Cells.Item(14,15).Value = [double]"0.2997993941754699"
$ws.Cells.Item(14,16).Value = [double]"0.29979939417547"
$ws.Cells.Item(14,17).Value = [double]"571.7114914731511"
$ws.Cells.Item(14,18).Value = [double]"5145.40342325836"
$ws.Cells.Item(14,19).Value = [double]"0.03168933636268883"
$ws.Cells.Item(14,20).Value = [double]"0.03168933636268884"

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Cxcl12"
$ws.Cells.Item(15,3).Value = "Itga4"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = [double]"3"
$ws.Cells.Item(15,6).Value = [double]"1"
$ws.Cells.Item(15,7).Value = [double]"24.16373066666667"
$ws.Cells.Item(15,8).Value = [double]"72.491192"
$ws.Cells.Item(15,9).Value = [double]"0.1057018025331343"
$ws.Cells.Item(15,10).Value = [double]"0.1057018025331344"
$ws.Cells.Item(15,11).Value = [double]"1"
$ws.Cells.Item(15,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(15,13).Value = [double]"0.07690566666666666"
$ws.Cells.Item(15,14).Value = [double]"0.230717"
$ws.Cells.Item(15,15).Value = [double]"0.0009744872400636476"
$ws.Cells.Item(15,16).Value = [double]"0.0009744872400636479"
$ws.Cells.Item(15,17).Value = [double]"1.858327816073778"
$ws.Cells.Item(15,18).Value = [double]"16.724950344664"
$ws.Cells.Item(15,19).Value = [double]"0.0001030050578202668"
$ws.Cells.Item(15,20).Value = [double]"0.0001030050578202668"

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Cxcl12"
$ws.Cells.Item(16,3).Value = "Itga4"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = [double]"3"
$ws.Cells.Item(16,6).Value = [double]"1"
$ws.Cells.Item(16,7).Value = [double]"24.16373066666667"
$ws.Cells.Item(16,8).Value = [double]"72.491192"
$ws.Cells.Item(16,9).Value = [double]"0.1057018025331343"
$ws.Cells.Item(16,10).Value = [double]"0.1057018025331344"
$ws.Cells.Item(16,11).Value = [double]"3"
$ws.Cells.Item(16,12).Value = [double]"1"
$ws.Cells.Item(16,13).Value = [double]"53.21452433333334"
$ws.Cells.Item(16,14).Value = [double]"159.643573"
$ws.Cells.Item(16,15).Value = [double]"0.6742919890890982"
$ws.Cells.Item(16,16).Value = [double]"0.6742919890890983"
$ws.Cells.Item(16,17).Value = [double]"1285.861433545446"
$ws.Cells.Item(16,18).Value = [double]"11572.75290190902"
$ws.Cells.Item(16,19).Value = [double]"0.07127387868037023"
$ws.Cells.Item(16,20).Value = [double]"0.07127387868037025"

# Row 17
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Cxcl12"
$ws.Cells.Item(17,3).Value = "Itga4"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = [double]"3"
$ws.Cells.Item(17,6).Value = [double]"1"
$ws.Cells.Item(17,7).Value = [double]"24.16373066666667"
$ws.Cells.Item(17,8).Value = [double]"72.491192"
$ws.Cells.Item(17,9).Value = [double]"0.1057018025331343"
$ws.Cells.Item(17,10).Value = [double]"0.1057018025331344"
$ws.Cells.Item(17,11).Value = [double]"3"
$ws.Cells.Item(17,12).Value = [double]"1"
$ws.Cells.Item(17,13).Value = [double]"1.967779333333333"
$ws.Cells.Item(17,14).Value = [double]"5.903338"
$ws.Cells.Item(17,15).Value = [double]"0.02493412949536815"
$ws.Cells.Item(17,16).Value = [double]"0.02493412949536816"
$ws.Cells.Item(17,17).Value = [double]"47.54888982209955"
$ws.Cells.Item(17,18).Value = [double]"427.940008398896"
$ws.Cells.Item(17,19).Value = [double]"0.002635582432255005"
$ws.Cells.Item(17,20).Value = [double]"0.002635582432255006"
